# OLX Monitor 2026-02-17 11:36 — append a fresh monitoring pass of listing
# rows (poqui / pokojewlublinie / dawnypatron) to the bottom of the first
# worksheet, reusing the prior pass (rows 7-14) as a formatting template so
# the styles (alignment + "stale listing" highlight) come along for free.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Duplicate the previous run's 8 detail rows (A7:H14) into A15:H22,
# carrying over cell styles/number formats exactly.
$src = $ws.Range("A7:H14")
$dst = $ws.Range("A15:H22")
$src.Copy($dst)

# Stamp the new pass with its own "last checked" timestamp.
$ws.Range("A15:A22").Value = "2026-02-17 11:36:19"
